$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$templateStart = 727
$startRow = 743
$name = "AHMED FOUZI AHMED ALHADDAD"
$certs = @("DSS1742","DSS1743","DSS1744","DSS1745","DSS1746","DSS1747","DSS1748","DSS1749")

# Copy the 8-row template block (same courses/dates layout, same "banded" style)
# down into the new block so formatting (styles, number formats) matches exactly.
$srcRange = $ws.Range("A" + $templateStart + ":E" + ($templateStart + 7))
$dstRange = $ws.Range("A" + $startRow + ":E" + ($startRow + 7))
$srcRange.Copy($dstRange)

# Fill in certificate numbers (column A) for all 8 rows first.
for ($i = 0; $i -lt 8; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $certs[$i]
}

# Then the new name (column B) for all 8 rows.
for ($i = 0; $i -lt 8; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $name
}

$ws.Range("A744").Select()
